$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.203.73'
$ws.Range('E2').Value = '  +1.86%  '
$ws.Range('D3').Value = '2.345.50'
$ws.Range('E3').Value = '  +5.95%  '
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').Value = '313.73'
$ws.Range('E5').Value = '  +6.42%  '
$ws.Range('D6').Value = '109.67'
$ws.Range('E6').Value = '  +1.74%  '
$ws.Range('E7').Value = '  +3.74%  '
$ws.Range('E8').Value = '  -0.20%  '
$ws.Range('D9').Value = '0.633'
$ws.Range('E9').Value = '  +6.25%  '
$ws.Range('D10').Value = '42.88'
$ws.Range('E10').Value = '  -1.71%  '
$ws.Range('D11').Value = '0.0937'
$ws.Range('E11').Value = '  +3.26%  '
$ws.Range('D12').Value = '8.87'
$ws.Range('E12').Value = '  +1.27%  '
$ws.Range('E13').Value = '  +9.33%  '
$ws.Range('E14').Value = '  +2.53%  '
$ws.Range('D15').Value = '16.26'
$ws.Range('E15').Value = '  +9.12%  '
$ws.Range('D16').Value = '2.701.15'
$ws.Range('E16').Value = '  +6.07%  '
$ws.Range('D17').Value = '2.342.86'
$ws.Range('E17').Value = '  +5.04%  '
$ws.Range('D18').Value = '43.170.82'
$ws.Range('E18').Value = '  +1.97%  '
$ws.Range('E19').Value = '  +3.96%  '
$ws.Range('D20').Value = '7.27'
$ws.Range('E20').Value = '  -1.34%  '
$ws.Range('D21').Value = '75.34'
$ws.Range('E21').Value = '  +3.82%  '
$ws.Range('D22').Value = '2.61'
$ws.Range('E22').Value = '  +14.80%  '
$ws.Range('D23').Value = '3.45'
$ws.Range('E23').Value = '  +0.16%  '
$ws.Range('D24').Value = '253.38'
$ws.Range('E24').Value = '  +11.34%  '
$ws.Range('D25').Value = '9.12'
$ws.Range('E25').Value = '  +1.41%  '
$ws.Range('D26').Value = '12.05'
$ws.Range('E26').Value = '  +4.26%  '
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('D28').Value = '39.38'
$ws.Range('E28').Value = '  +2.34%  '
$ws.Range('E29').Value = '  +1.64%  '
$ws.Range('D30').Value = '22.40'
$ws.Range('D31').Value = '174.68'
$ws.Range('E32').Value = '  -0.97%  '
$ws.Range('D33').Value = '0.0929'
$ws.Range('E33').Value = '  +5.15%  '
$ws.Range('E34').Value = '  +8.91%  '
$ws.Range('D35').Value = '0.132'
$ws.Range('E35').Value = '  +5.75%  '
$ws.Range('E36').Value = '  -0.41%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').Value = '4.16'
$ws.Range('E37').Value = '  -4.25%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.0377'
$ws.Range('E38').Value = '  +4.20%  '
$ws.Range('E39').Value = '  +1.70%  '
$ws.Range('E40').Value = '  +10.64%  '
$ws.Range('D41').Value = '73.04'
$ws.Range('E41').Value = '  +3.07%  '
$ws.Range('D42').Value = '1.47'
$ws.Range('E42').Value = '  +13.40%  '
$ws.Range('E43').Value = '  +1.54%  '
$ws.Range('D44').Value = '12.90'
$ws.Range('E44').Value = '  +2.10%  '
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('D46').Value = '5.64'
$ws.Range('E46').Value = '  +4.25%  '
$ws.Range('D47').Value = '9.31'
$ws.Range('E47').Value = '  +10.20%  '
$ws.Range('D48').Value = '110.95'
$ws.Range('E48').Value = '  +7.67%  '
$ws.Range('E49').Value = '  -0.66%  '
$ws.Range('E50').Value = '  +4.00%  '
$ws.Range('D51').Value = '69.96'
$ws.Range('E51').Value = '  +4.96%  '
